$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    4   = -21.05549999999998
    6   = -22.71810000000001
    7   = -20.21649999999997
    16  = -21.9523
    20  = -19.92889999999999
    28  = -21.92770000000001
    29  = -21.29709999999997
    32  = -21.1996
    40  = -20.3163
    46  = -21.7836
    51  = -21.61679999999999
    52  = -22.18260000000001
    57  = -22.68830000000002
    59  = -22.13810000000001
    62  = -22.09550000000001
    66  = -21.4546
    73  = -20.26939999999999
    74  = -21.64319999999998
    92  = -21.60629999999999
    100 = -22.08800000000002
}

foreach ($row in $changes.Keys) {
    $ws.Cells.Item($row, 1).Value = $changes[$row]
}
